# Append two new observation rows (5 and 6) to the "Artfynd" sheet,
# matching the rows already present (Tretåig hackspett / Picoides tridactylus).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y and AA hold dates stored as plain text (e.g. "2026-02-23"),
# not real Excel date serials, so force text formatting before writing
# them to stop Excel from auto-converting the strings to date numbers.
$ws.Range("Y5:Y6").NumberFormat = "@"
$ws.Range("AA5:AA6").NumberFormat = "@"

# ---- Row 5 ----
$ws.Range("A5").Value = 131273746
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("I5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "färska spår"
$ws.Range("N5").Value = ""
$ws.Range("P5").Value = "Sims fäbodar, Dlr"
$ws.Range("Q5").Value = 515255
$ws.Range("R5").Value = 6705041
$ws.Range("S5").Value = 50
$ws.Range("T5").Value = "Dalarna"
$ws.Range("U5").Value = "Borlänge"
$ws.Range("V5").Value = "Dalarna"
$ws.Range("W5").Value = "Stora Tuna"
$ws.Range("Y5").Value = "2026-02-23"
$ws.Range("AA5").Value = "2026-02-23"
$ws.Range("AC5").Value = "Ringhack på tall."
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = "Anna-Lena Thommson"
$ws.Range("AX5").Value = "Anna-Lena Thommson"
$ws.Range("AY5").Value = ""

# ---- Row 6 ----
$ws.Range("A6").Value = 131273722
$ws.Range("B6").Value = 57884
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("I6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "äldre spår"
$ws.Range("N6").Value = ""
$ws.Range("P6").Value = "Sims fäbodar, Dlr"
$ws.Range("Q6").Value = 515365
$ws.Range("R6").Value = 6705054
$ws.Range("S6").Value = 50
$ws.Range("T6").Value = "Dalarna"
$ws.Range("U6").Value = "Borlänge"
$ws.Range("V6").Value = "Dalarna"
$ws.Range("W6").Value = "Stora Tuna"
$ws.Range("Y6").Value = "2026-02-23"
$ws.Range("AA6").Value = "2026-02-23"
$ws.Range("AC6").Value = "Ringhack på tall."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Anna-Lena Thommson"
$ws.Range("AX6").Value = "Anna-Lena Thommson"
$ws.Range("AY6").Value = ""

Write-Output "Added rows 5 and 6 to $($ws.Name)"
